# Generate Report for Handoff
# Updates the localization-status report with a new handoff id
# (f75a1d65-2938-47d3-af59-1ca3787b8220 -> 0cbbba66-5284-457a-8025-6fc096adc98d)
# and refreshed handoff timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldId = "f75a1d65-2938-47d3-af59-1ca3787b8220"
$newId = "0cbbba66-5284-457a-8025-6fc096adc98d"

$oldHash = "d21e233fe8288ec355c6e6808a0d52f27cbca0c6"
$newHash = "43b760d6203a45b8224c2b31b13f5a71b18e4eb0"

# ---------------------------------------------------------------------------
# Overview sheet: file name + latest handoff date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$mdAddr = "https://github.com/OpenLocalizationTest/oltest/blob/bf98ce580bbe82a06712c39c76dcdc933a25a20d/e2e/$oldId.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddr, "", "", "$newId.md")

$wsOverview.Range("D2").Value = "2016-59-13 10:59:07"

# ---------------------------------------------------------------------------
# zh-cn sheet: file name, xlf target file name, handoff datetime
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c3ce99c6783852cd367eea4abd376a5a4c4a1a68/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldHash.zh-cn.xlf"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddr, "", "", "$newId.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $mdAddr, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhXlfAddr, "", "", "$newId.$newHash.zh-cn.xlf")

$wsZhCn.Range("E2").Value = "2016-03-13 10:59:03"

# ---------------------------------------------------------------------------
# de-de sheet: file name, xlf target file name, handoff datetime
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a72f04208a034c340f3b286cca2be3f9a2278fcd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldHash.de-de.xlf"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddr, "", "", "$newId.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $mdAddr, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deXlfAddr, "", "", "$newId.$newHash.de-de.xlf")

$wsDeDe.Range("E2").Value = "2016-03-13 10:59:07"
